$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '62.688.98'
Set-TextValue $ws.Range('E2') '  +1.85%  '
Set-TextValue $ws.Range('D3') '3.460.61'
Set-TextValue $ws.Range('E3') '  +2.05%  '
Set-TextValue $ws.Range('E4') '  +0.04%  '
Set-TextValue $ws.Range('D5') '578.15'
Set-TextValue $ws.Range('E5') '  +0.44%  '
Set-TextValue $ws.Range('D6') '146.21'
Set-TextValue $ws.Range('E6') '  +3.59%  '
Set-TextValue $ws.Range('E7') '  -0.06%  '
Set-TextValue $ws.Range('D8') '0.482'
Set-TextValue $ws.Range('E8') '  +2.05%  '
Set-TextValue $ws.Range('D9') '7.61'
Set-TextValue $ws.Range('E9') '  -0.66%  '
Set-TextValue $ws.Range('E10') '  +1.78%  '
Set-TextValue $ws.Range('D11') '0.399'
Set-TextValue $ws.Range('E11') '  +3.59%  '
Set-TextValue $ws.Range('D12') '4.051.32'
Set-TextValue $ws.Range('E12') '  +2.07%  '
Set-TextValue $ws.Range('D13') '29.76'
Set-TextValue $ws.Range('E13') '  +4.75%  '
Set-TextValue $ws.Range('E14') '  +2.38%  '
Set-TextValue $ws.Range('D15') '3.461.92'
Set-TextValue $ws.Range('E15') '  +1.91%  '
Set-TextValue $ws.Range('D16') '0.0000170'
Set-TextValue $ws.Range('E16') '  +0.27%  '
Set-TextValue $ws.Range('D17') '62.775.40'
Set-TextValue $ws.Range('E17') '  +1.91%  '
Set-TextValue $ws.Range('D18') '6.35'
Set-TextValue $ws.Range('E18') '  +3.67%  '
Set-TextValue $ws.Range('D19') '14.39'
Set-TextValue $ws.Range('E19') '  +5.70%  '
Set-TextValue $ws.Range('D20') '9.22'
Set-TextValue $ws.Range('E20') '  +2.56%  '
Set-TextValue $ws.Range('D21') '388.38'
Set-TextValue $ws.Range('E21') '  -0.72%  '
Set-TextValue $ws.Range('D22') '0.563'
Set-TextValue $ws.Range('E22') '  +2.44%  '
Set-TextValue $ws.Range('D23') '74.95'
Set-TextValue $ws.Range('E23') '  -0.06%  '
Set-TextValue $ws.Range('D25') '3.607.90'
Set-TextValue $ws.Range('E25') '  +2.32%  '
Set-TextValue $ws.Range('E26') '  +1.18%  '
Set-TextValue $ws.Range('D27') '0.178'
Set-TextValue $ws.Range('D28') '7.59'
Set-TextValue $ws.Range('E28') '  +4.39%  '
Set-TextValue $ws.Range('E29') '  +0.07%  '
Set-TextValue $ws.Range('D30') '8.11'
Set-TextValue $ws.Range('E30') '  +0.94%  '
Set-TextValue $ws.Range('E32') '  -0.01%  '
Set-TextValue $ws.Range('E33') '  -0.10%  '
Set-TextValue $ws.Range('D34') '23.74'
Set-TextValue $ws.Range('E34') '  +1.92%  '
Set-TextValue $ws.Range('D35') '7.08'
Set-TextValue $ws.Range('E35') '  +2.64%  '
Set-TextValue $ws.Range('E36') '  +5.09%  '
Set-TextValue $ws.Range('E37') '  +6.79%  '
Set-TextValue $ws.Range('D38') '31.33'
Set-TextValue $ws.Range('E38') '  +20.77%  '
Set-TextValue $ws.Range('D39') '169.79'
Set-TextValue $ws.Range('E39') '  +0.67%  '
Set-TextValue $ws.Range('D40') '3.499.49'
Set-TextValue $ws.Range('E40') '  +2.18%  '
Set-TextValue $ws.Range('D41') '0.0768'
Set-TextValue $ws.Range('E41') '  +0.46%  '
Set-TextValue $ws.Range('D42') '0.797'
Set-TextValue $ws.Range('D43') '4.48'
Set-TextValue $ws.Range('E43') '  +1.59%  '
Set-TextValue $ws.Range('D44') '42.14'
Set-TextValue $ws.Range('E44') '  -0.80%  '
Set-TextValue $ws.Range('E45') '  +3.37%  '
Set-TextValue $ws.Range('E46') '  +1.97%  '
Set-TextValue $ws.Range('D47') '2.594.08'
Set-TextValue $ws.Range('E47') '  +4.57%  '
Set-TextValue $ws.Range('D48') '23.36'
Set-TextValue $ws.Range('E48') '  +2.21%  '
Set-TextValue $ws.Range('E49') '  +1.62%  '
Set-TextValue $ws.Range('D50') '2.19'
Set-TextValue $ws.Range('E50') '  +8.18%  '
Set-TextValue $ws.Range('E51') '  +0.01%  '
